# Apply the 2022 column (S) addition and updated 2020/2021 (Q/R) figures
# to the Consumer Price Index worksheet, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing Q (2020) and R (2021) values for rows 5-14 ---

$qrUpdates = @(
    @{ Row = 5;  Q = 117.60684979252385; R = 113.34848864817617 },
    @{ Row = 6;  Q = 114.77319768114526; R = 115.06069350712495 },
    @{ Row = 7;  Q = 116.40044011407315; R = 114.29658549692938 },
    @{ Row = 8;  Q = 117.53828537152096; R = 113.75761785228545 },
    @{ Row = 9;  Q = 117.42206669681742; R = 113.98264089946031 },
    @{ Row = 10; Q = 113.98326995089161; R = 113.92720567782911 },
    @{ Row = 11; Q = 123.488978736909;   R = 114.17226706705155 },
    @{ Row = 12; Q = 118.12340252754679; R = 114.45153946490467 },
    @{ Row = 13; Q = 118.87059844457349; R = 112.69493421065988 },
    @{ Row = 14; Q = 114.06377070452145; R = 113.95067699644588 }
)

foreach ($u in $qrUpdates) {
    $ws.Cells.Item($u.Row, 17).Value() = $u.Q   # column Q = 17
    $ws.Cells.Item($u.Row, 18).Value() = $u.R   # column R = 18
}

# --- 2. Add the new S column (2022) ---

# Header cell S4: copy formatting from R4 (year header style) then set value
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S4").Value() = 2022

# Data cells S5:S13: copy formatting from the default column style (column N, style 4)
$sDataUpdates = @(
    @{ Row = 5;  Value = 115.8 },
    @{ Row = 6;  Value = 115.2 },
    @{ Row = 7;  Value = 115.4 },
    @{ Row = 8;  Value = 111.8 },
    @{ Row = 9;  Value = 116.8 },
    @{ Row = 10; Value = 108.2 },
    @{ Row = 11; Value = 111 },
    @{ Row = 12; Value = 115.8 },
    @{ Row = 13; Value = 117.9 }
)

foreach ($u in $sDataUpdates) {
    $ws.Cells.Item($u.Row, 19).Value() = $u.Value   # column S = 19
}

# S14: copy formatting from R14 (bottom border / thick bottom row style) then set value
$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S14").Value() = 112.4

# --- 3. Update the active cell selection (T6 -> T4) ---
$null = $ws.Range("T4").Select()
